$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 112-171: shift dates (and a few associated value changes)
$ws.Cells.Item(112,4).Value = 44839   # D112
$ws.Cells.Item(112,11).Value = 700   # K112
$ws.Cells.Item(112,12).Value = 800   # L112
$ws.Cells.Item(112,13).Value = 750   # M112
$ws.Cells.Item(112,16).Value = 750   # P112
$ws.Cells.Item(113,4).Value = 44839   # D113
$ws.Cells.Item(113,11).Value = 600   # K113
$ws.Cells.Item(113,12).Value = 600   # L113
$ws.Cells.Item(113,13).Value = 600   # M113
$ws.Cells.Item(113,16).Value = 600   # P113
$ws.Cells.Item(114,4).Value = 44358   # D114
$ws.Cells.Item(115,4).Value = 44358   # D115
$ws.Cells.Item(116,4).Value = 44327   # D116
$ws.Cells.Item(117,4).Value = 44327   # D117
$ws.Cells.Item(118,4).Value = 44460   # D118
$ws.Cells.Item(119,4).Value = 44460   # D119
$ws.Cells.Item(120,4).Value = 44160   # D120
$ws.Cells.Item(121,4).Value = 44160   # D121
$ws.Cells.Item(122,4).Value = 44609   # D122
$ws.Cells.Item(122,11).Value = 600   # K122
$ws.Cells.Item(122,12).Value = 700   # L122
$ws.Cells.Item(122,13).Value = 650   # M122
$ws.Cells.Item(122,16).Value = 650   # P122
$ws.Cells.Item(123,4).Value = 44609   # D123
$ws.Cells.Item(123,11).Value = 500   # K123
$ws.Cells.Item(123,12).Value = 500   # L123
$ws.Cells.Item(123,13).Value = 500   # M123
$ws.Cells.Item(123,16).Value = 500   # P123
$ws.Cells.Item(124,4).Value = 44771   # D124
$ws.Cells.Item(124,11).Value = 700   # K124
$ws.Cells.Item(124,12).Value = 800   # L124
$ws.Cells.Item(124,13).Value = 750   # M124
$ws.Cells.Item(124,16).Value = 750   # P124
$ws.Cells.Item(125,4).Value = 44771   # D125
$ws.Cells.Item(125,11).Value = 600   # K125
$ws.Cells.Item(125,12).Value = 600   # L125
$ws.Cells.Item(125,13).Value = 600   # M125
$ws.Cells.Item(125,16).Value = 600   # P125
$ws.Cells.Item(126,4).Value = 44308   # D126
$ws.Cells.Item(127,4).Value = 44308   # D127
$ws.Cells.Item(128,4).Value = 44224   # D128
$ws.Cells.Item(129,4).Value = 44224   # D129
$ws.Cells.Item(130,4).Value = 44166   # D130
$ws.Cells.Item(131,4).Value = 44166   # D131
$ws.Cells.Item(132,4).Value = 44435   # D132
$ws.Cells.Item(132,10).Value = 200   # J132
$ws.Cells.Item(133,4).Value = 44435   # D133
$ws.Cells.Item(133,10).Value = 100   # J133
$ws.Cells.Item(134,4).Value = 44442   # D134
$ws.Cells.Item(134,10).Value = 300   # J134
$ws.Cells.Item(135,4).Value = 44442   # D135
$ws.Cells.Item(135,10).Value = 150   # J135
$ws.Cells.Item(136,4).Value = 44336   # D136
$ws.Cells.Item(137,4).Value = 44336   # D137
$ws.Cells.Item(138,4).Value = 44252   # D138
$ws.Cells.Item(139,4).Value = 44252   # D139
$ws.Cells.Item(140,4).Value = 44694   # D140
$ws.Cells.Item(141,4).Value = 44694   # D141
$ws.Cells.Item(142,4).Value = 44405   # D142
$ws.Cells.Item(143,4).Value = 44405   # D143
$ws.Cells.Item(144,4).Value = 44679   # D144
$ws.Cells.Item(145,4).Value = 44679   # D145
$ws.Cells.Item(146,4).Value = 44231   # D146
$ws.Cells.Item(147,4).Value = 44231   # D147
$ws.Cells.Item(148,4).Value = 44334   # D148
$ws.Cells.Item(149,4).Value = 44334   # D149
$ws.Cells.Item(150,4).Value = 44194   # D150
$ws.Cells.Item(151,4).Value = 44194   # D151
$ws.Cells.Item(152,4).Value = 44330   # D152
$ws.Cells.Item(153,4).Value = 44330   # D153
$ws.Cells.Item(154,4).Value = 44274   # D154
$ws.Cells.Item(155,4).Value = 44274   # D155
$ws.Cells.Item(156,4).Value = 44391   # D156
$ws.Cells.Item(157,4).Value = 44391   # D157
$ws.Cells.Item(158,4).Value = 44433   # D158
$ws.Cells.Item(159,4).Value = 44433   # D159
$ws.Cells.Item(160,4).Value = 44203   # D160
$ws.Cells.Item(161,4).Value = 44203   # D161
$ws.Cells.Item(162,4).Value = 44355   # D162
$ws.Cells.Item(162,15).Value = 'Región de Ñuble'   # O162
$ws.Cells.Item(163,4).Value = 44355   # D163
$ws.Cells.Item(163,15).Value = 'Región de Ñuble'   # O163
$ws.Cells.Item(164,4).Value = 44565   # D164
$ws.Cells.Item(164,15).Value = 'Región Metropolitana'   # O164
$ws.Cells.Item(165,4).Value = 44565   # D165
$ws.Cells.Item(165,15).Value = 'Región Metropolitana'   # O165
$ws.Cells.Item(166,4).Value = 44187   # D166
$ws.Cells.Item(167,4).Value = 44187   # D167
$ws.Cells.Item(168,4).Value = 44553   # D168
$ws.Cells.Item(169,4).Value = 44553   # D169
$ws.Cells.Item(170,4).Value = 44292   # D170
$ws.Cells.Item(171,4).Value = 44292   # D171

# Append new rows 172-173 with full data
$ws.Cells.Item(172,1).Value = 11   # A172
$ws.Cells.Item(172,2).Value = 'Vega Monumental Concepción'   # B172
$ws.Cells.Item(172,3).Value = 'Bíobío'   # C172
$ws.Cells.Item(172,4).Value = 44453   # D172
$ws.Cells.Item(172,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(172,5).Value = 8   # E172
$ws.Cells.Item(172,6).Value = 100112044   # F172
$ws.Cells.Item(172,7).Value = 'Perejil'   # G172
$ws.Cells.Item(172,8).Value = 'Sin especificar'   # H172
$ws.Cells.Item(172,9).Value = 'Primera'   # I172
$ws.Cells.Item(172,10).Value = 200   # J172
$ws.Cells.Item(172,11).Value = 600   # K172
$ws.Cells.Item(172,12).Value = 700   # L172
$ws.Cells.Item(172,13).Value = 650   # M172
$ws.Cells.Item(172,14).Value = '$/atado 0,5 a 1 kilo'   # N172
$ws.Cells.Item(172,15).Value = 'Región de Ñuble'   # O172
$ws.Cells.Item(172,16).Value = 650   # P172
$ws.Cells.Item(172,17).Value = 1   # Q172
$ws.Cells.Item(172,18).Value = 'Hortaliza'   # R172
$ws.Cells.Item(173,1).Value = 11   # A173
$ws.Cells.Item(173,2).Value = 'Vega Monumental Concepción'   # B173
$ws.Cells.Item(173,3).Value = 'Bíobío'   # C173
$ws.Cells.Item(173,4).Value = 44453   # D173
$ws.Cells.Item(173,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(173,5).Value = 8   # E173
$ws.Cells.Item(173,6).Value = 100112044   # F173
$ws.Cells.Item(173,7).Value = 'Perejil'   # G173
$ws.Cells.Item(173,8).Value = 'Sin especificar'   # H173
$ws.Cells.Item(173,9).Value = 'Segunda'   # I173
$ws.Cells.Item(173,10).Value = 100   # J173
$ws.Cells.Item(173,11).Value = 500   # K173
$ws.Cells.Item(173,12).Value = 500   # L173
$ws.Cells.Item(173,13).Value = 500   # M173
$ws.Cells.Item(173,14).Value = '$/atado 0,5 a 1 kilo'   # N173
$ws.Cells.Item(173,15).Value = 'Región de Ñuble'   # O173
$ws.Cells.Item(173,16).Value = 500   # P173
$ws.Cells.Item(173,17).Value = 1   # Q173
$ws.Cells.Item(173,18).Value = 'Hortaliza'   # R173
